# Update the K column (column G) values for rows 2-13.
# These values represent recalculated "K" (previously Strike#) data
# regenerated as part of refreshing std/mean and s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 1
    4  = 4
    5  = 4
    6  = 5
    7  = 4
    8  = 2
    9  = 2
    10 = 1
    11 = 4
    12 = 1
    13 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
